$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A62:A83").ClearContents()
$ws.Range("C62:D83").ClearContents()
